$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove trailing Russian age-unit word ("года"/"лет") from the age-range
# labels in column A, rows 2-22, leaving just the numeric range.
$labels = @{
    2  = "0 - 4"
    3  = "5 - 9"
    4  = "10 - 14"
    5  = "15 - 19"
    6  = "20 - 24"
    7  = "25 - 29"
    8  = "30 - 34"
    9  = "35 - 39"
    10 = "40 - 44"
    11 = "45 - 49"
    12 = "50 - 54"
    13 = "55 - 59"
    14 = "60 - 64"
    15 = "65 - 69"
    16 = "70 - 74"
    17 = "75 - 79"
    18 = "80 - 84"
    19 = "85 - 89"
    20 = "90 - 94"
    21 = "95 - 99"
    22 = "100+"
}

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 1).Value = $labels[$row]
}
